# Add first after-work walk of the year.
#
# The workbook tracks a monthly walking log on "Sheet1":
#   column A = month name, column B/C = days/walks-per-week inputs,
#   column D/E = computed running totals feeding the "total distance" line
#   on the Chart1 chart, and column G holds the distance of the (first)
#   after-work walk logged for that month, mirrored into column F
#   (F1 = G1, F2 = F1 + G2, ...) which feeds the second line series on
#   the same chart.
#
# Recording the year's first after-work walk (26.5) for January means
# writing it into Sheet1!G1 - the F-column mirror and the chart series
# recalculate from that single input.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G1").Value = 26.5

# Nudge the embedded chart on the "Chart1" sheet to pick up the new value
# (best-effort; harmless if the host doesn't need/implement it).
try {
    $chartSheet = $wb.Worksheets.Item("Chart1")
    $chartObjects = $chartSheet.ChartObjects()
    if ($chartObjects.Count -gt 0) {
        $chartObjects.Item(1).Chart.Refresh()
    }
} catch {
    # Refresh isn't essential to the edit - ignore if unavailable.
}
